$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply per-cell value updates as described by the source diff.
# Column D holds price strings that must stay text; force text format
# before assigning any value that Excel would otherwise auto-convert
# to a number (losing trailing zeros / exact formatting).
$ws.Range("D2").Value = '57.958.10'
$ws.Range("E2").Value = '  -1.61%  '
$ws.Range("D3").Value = '2.449.47'
$ws.Range("E3").Value = '  -3.81%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.05'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.72'
$ws.Range("E6").Value = '  -3.63%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.565'
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0975'
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("E10").Value = '  -2.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.97'
$ws.Range("E11").Value = '  -4.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.322'
$ws.Range("E12").Value = '  -4.22%  '
$ws.Range("D13").Value = '2.882.87'
$ws.Range("E13").Value = '  -3.76%  '
$ws.Range("D14").Value = '57.900.01'
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.56'
$ws.Range("E15").Value = '  -3.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000132'
$ws.Range("E16").Value = '  -3.19%  '
$ws.Range("D17").Value = '2.452.49'
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.37'
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.11'
$ws.Range("E19").Value = '  -2.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '313.67'
$ws.Range("E20").Value = '  -3.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.14'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.09'
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.405'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("D26").Value = '2.563.10'
$ws.Range("E26").Value = '  -3.18%  '
$ws.Range("E27").Value = '  -2.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.24'
$ws.Range("E28").Value = '  -3.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '174.93'
$ws.Range("E29").Value = '  +3.79%  '
$ws.Range("D30").Value = '0.0₃0734'
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("E31").Value = '  -2.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.17'
$ws.Range("E32").Value = '  -3.61%  '
$ws.Range("E33").Value = '  -7.34%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.81'
$ws.Range("E36").Value = '  -3.05%  '
$ws.Range("E37").Value = '  -7.78%  '
$ws.Range("E38").Value = '  -5.39%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("E40").Value = '  +2.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.44'
$ws.Range("E41").Value = '  -4.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.38'
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '125.92'
$ws.Range("E43").Value = '  -4.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.584'
$ws.Range("E44").Value = '  -3.64%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.79'
$ws.Range("E45").Value = '  -6.11%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '257.66'
$ws.Range("E46").Value = '  -9.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0921'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("E49").Value = '  -3.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.01'
$ws.Range("E50").Value = '  -5.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.31'
